$wb = $excel.ActiveWorkbook

# Overview sheet: row 3 is the 71c0c1a2-... file, update Status (B and C) to "Handed back: in sync with en-US"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"

# zh-cn sheet: update status for row 3, and refresh handback datetimes in column G for rows 2 and 3
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("B3").Value = "Handed back: in sync with en-US"
$wsZh.Range("G2").Value = "2016-02-23 09:31:14"
$wsZh.Range("G3").Value = "2016-02-23 09:31:14"

# de-de sheet: update status for row 3, and refresh handback datetimes in column G for rows 2 and 3
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("B3").Value = "Handed back: in sync with en-US"
$wsDe.Range("G2").Value = "2016-02-23 09:31:46"
$wsDe.Range("G3").Value = "2016-02-23 09:31:46"
